$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C21").Value = '音频'
$ws.Range("C22").Value = '操作'
$ws.Range("C24").Value = '全屏'
$ws.Range("C25").Value = '无边框'
$ws.Range("A24").Value = 'Fullscreen'
$ws.Range("B24").Value = 'Fullscreen'
$ws.Range("A20").Value = 'Video'
$ws.Range("B20").Value = 'Video'
$ws.Range("A21").Value = 'Audio'
$ws.Range("B21").Value = 'Audio'
$ws.Range("A22").Value = 'Controls'
$ws.Range("B22").Value = 'Controls'
$ws.Range("A25").Value = 'Borderless'
$ws.Range("B25").Value = 'Borderless'
$ws.Range("A26").Value = 'Disabled'
$ws.Range("B26").Value = 'Disabled'
$ws.Range("A27").Value = 'Enabled'
$ws.Range("B27").Value = 'Enabled'
$ws.Range("A28").Value = 'Adaptive'
$ws.Range("B28").Value = 'Adaptive'
$ws.Range("A29").Value = 'Nvidia DLSS'
$ws.Range("B29").Value = 'Nvidia DLSS'
$ws.Range("C29").Value = 'Nvidia DLSS'
$ws.Range("C28").Value = '自适应'
$ws.Range("C26").Value = '关闭'
$ws.Range("C27").Value = '开启'
$ws.Range("A23").Value = 'Return'
$ws.Range("B23").Value = 'Return'
$ws.Range("C23").Value = '返回'
$ws.Range("C20").Value = '图像'
$ws.Range("B30").Value = 'Master'
$ws.Range("A30").Value = 'Master_Music'
$ws.Range("C30").Value = '主音量'
$ws.Range("A31").Value = 'Music'
$ws.Range("B31").Value = 'Music'
$ws.Range("C31").Value = '音乐'
$ws.Range("A32").Value = 'SFX'
$ws.Range("B32").Value = 'SFX'
$ws.Range("C32").Value = '音效'
$ws.Range("A33").Value = 'Game'
$ws.Range("B33").Value = 'Game'
$ws.Range("C33").Value = '游戏'
$ws.Range("A34").Value = 'Paused_Settings'
$ws.Range("B34").Value = 'settings'
$ws.Range("C34").Value = '设置'
$ws.Range("A35").Value = 'Resume'
$ws.Range("B35").Value = 'Resume'
$ws.Range("C35").Value = '恢复游戏'
$ws.Range("A36").Value = 'Restart'
$ws.Range("B36").Value = 'Restart'
$ws.Range("C36").Value = '重新开始'
$ws.Range("A37").Value = 'Main Menu'
$ws.Range("B37").Value = 'Main Menu'
$ws.Range("C37").Value = '主菜单'
$ws.Range("A38").Value = 'Game Paused'
$ws.Range("B38").Value = 'Game Paused'
$ws.Range("C38").Value = '游戏暂停'

$ws.Range("C38").Select() | Out-Null
